$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.315.50"
$ws.Range("E2").Value = "  +1.48%  "
$ws.Range("D3").Value = "1.622.47"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.488"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("E8").Value = "  +1.61%  "
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.90%  "
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").Value = "1.848.14"
$ws.Range("E12").Value = "  +1.86%  "
$ws.Range("D13").Value = "1.621.15"
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.519"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").Value = "26.320.21"
$ws.Range("E16").Value = "  +1.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.36%  "
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "203.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("E21").Value = "  +1.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.52%  "
$ws.Range("E23").Value = "  +1.14%  "
$ws.Range("E24").Value = "  +6.47%  "
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("E29").Value = "  +2.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0531"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +12.00%  "
$ws.Range("E31").Value = "  +0.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("E34").Value = "  +2.45%  "
$ws.Range("E35").Value = "  +2.20%  "
$ws.Range("D36").Value = "1.181.65"
$ws.Range("E36").Value = "  +5.23%  "
$ws.Range("E37").Value = "  +1.85%  "
$ws.Range("E38").Value = "  +3.51%  "
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.496"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.788"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.43%  "
$ws.Range("E43").Value = "  +5.23%  "
$ws.Range("D44").Value = "1.758.77"
$ws.Range("E44").Value = "  +1.87%  "
$ws.Range("E45").Value = "  +1.29%  "
$ws.Range("E46").Value = "  +15.54%  "
$ws.Range("E47").Value = "  +2.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.75%  "
$ws.Range("E49").Value = "  +1.09%  "
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("E51").Value = "  -0.38%  "
